# Auto-generated script applying numeric corrections to the Leve profit tables
# across all 8 sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 256688.56
$ws.Range("I6").Value = 276274.75
$ws.Range("K6").Value = 828824.25
$ws.Range("M6").Value = -828712.25
$ws.Range("H8").Value = 26.8
$ws.Range("I8").Value = 27.555555
$ws.Range("J8").Value = 20
$ws.Range("K8").Value = 82.66666499999999
$ws.Range("L8").Value = 60
$ws.Range("M8").Value = 56.33333500000001
$ws.Range("N8").Value = -338
$ws.Range("H12").Value = 25102.5
$ws.Range("I12").Value = 25102.5
$ws.Range("K12").Value = 25102.5
$ws.Range("M12").Value = -24932.5
$ws.Range("H17").Value = 2033.9231
$ws.Range("J17").Value = 4100
$ws.Range("L17").Value = 12300
$ws.Range("N17").Value = -12636
$ws.Range("H28").Value = 477.4
$ws.Range("I28").Value = 521
$ws.Range("K28").Value = 521
$ws.Range("M28").Value = -36
$ws.Range("H33").Value = 98.0625
$ws.Range("I33").Value = 98.0625
$ws.Range("K33").Value = 98.0625
$ws.Range("M33").Value = 130.9375
$ws.Range("H39").Value = 243.5
$ws.Range("I39").Value = 158
$ws.Range("J39").Value = 500
$ws.Range("K39").Value = 474
$ws.Range("L39").Value = 1500
$ws.Range("M39").Value = -178
$ws.Range("N39").Value = -2092
$ws.Range("H40").Value = 1510.7407
$ws.Range("I40").Value = 1324.5834
$ws.Range("K40").Value = 1324.5834
$ws.Range("M40").Value = -1149.5834
$ws.Range("H51").Value = 4124.982
$ws.Range("H62").Value = 8535.875
$ws.Range("I62").Value = 10448.833
$ws.Range("J62").Value = 2797
$ws.Range("K62").Value = 10448.833
$ws.Range("L62").Value = 2797
$ws.Range("M62").Value = -9824.833000000001
$ws.Range("N62").Value = -4045
$ws.Range("H65").Value = 8535.875
$ws.Range("I65").Value = 10448.833
$ws.Range("J65").Value = 2797
$ws.Range("K65").Value = 52244.165
$ws.Range("L65").Value = 13985
$ws.Range("M65").Value = -49124.165
$ws.Range("N65").Value = -20225
$ws.Range("H86").Value = 7279.4
$ws.Range("J86").Value = 4450
$ws.Range("L86").Value = 4450
$ws.Range("N86").Value = -6696
$ws.Range("H89").Value = 7279.4
$ws.Range("J89").Value = 4450
$ws.Range("L89").Value = 22250
$ws.Range("N89").Value = -33482
$ws.Range("H92").Value = 1085.1578
$ws.Range("I92").Value = 319
$ws.Range("K92").Value = 319
$ws.Range("M92").Value = 929
$ws.Range("H98").Value = 2317.6428
$ws.Range("J98").Value = 1154
$ws.Range("L98").Value = 1154
$ws.Range("N98").Value = -4150
$ws.Range("H107").Value = 8382.723
$ws.Range("I107").Value = 7259.933
$ws.Range("J107").Value = 13996.667
$ws.Range("K107").Value = 7259.933
$ws.Range("L107").Value = 13996.667
$ws.Range("M107").Value = -5339.933
$ws.Range("N107").Value = -17836.667
$ws.Range("H111").Value = 1720.2
$ws.Range("I111").Value = 1642.25
$ws.Range("K111").Value = 4926.75
$ws.Range("M111").Value = -1859.75
$ws.Range("H116").Value = 13022.862
$ws.Range("I116").Value = 3568.3076
$ws.Range("J116").Value = 20704.688
$ws.Range("K116").Value = 3568.3076
$ws.Range("L116").Value = 20704.688
$ws.Range("M116").Value = -126.3076000000001
$ws.Range("N116").Value = -27588.688
$ws.Range("H121").Value = 2765.1
$ws.Range("J121").Value = 3579.1428
$ws.Range("L121").Value = 10737.4284
$ws.Range("N121").Value = -14231.4284
$ws.Range("H122").Value = 2317.6428
$ws.Range("J122").Value = 1154
$ws.Range("L122").Value = 3462
$ws.Range("N122").Value = -8362
$ws.Range("H125").Value = 3470.7273
$ws.Range("I125").Value = 2677.8
$ws.Range("K125").Value = 24100.2
$ws.Range("M125").Value = -21640.2
$ws.Range("H127").Value = 47313.22
$ws.Range("I127").Value = 49395.637
$ws.Range("K127").Value = 148186.911
$ws.Range("M127").Value = -143226.911
$ws.Range("H131").Value = 27125.125
$ws.Range("I131").Value = 2400.2
$ws.Range("K131").Value = 7200.599999999999
$ws.Range("M131").Value = -2160.599999999999
$ws.Range("H132").Value = 6252.8667
$ws.Range("I132").Value = 6976.0264
$ws.Range("K132").Value = 20928.0792
$ws.Range("M132").Value = -18398.0792
$ws.Range("H137").Value = 20005250
$ws.Range("I137").Value = 55557150
$ws.Range("J137").Value = 7303.8438
$ws.Range("K137").Value = 166671450
$ws.Range("L137").Value = 21911.5314
$ws.Range("M137").Value = -166668900
$ws.Range("N137").Value = -27011.5314
$ws.Range("H138").Value = 4565.5293
$ws.Range("J138").Value = 4394.5938
$ws.Range("L138").Value = 13183.7814
$ws.Range("N138").Value = -23463.7814
$ws.Range("H141").Value = 14712.143
$ws.Range("I141").Value = 14328.333
$ws.Range("K141").Value = 42984.999
$ws.Range("M141").Value = -37804.999

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1146.6666
$ws.Range("I2").Value = 970
$ws.Range("J2").Value = 1500
$ws.Range("K2").Value = 970
$ws.Range("L2").Value = 1500
$ws.Range("M2").Value = -857
$ws.Range("N2").Value = -1726
$ws.Range("H32").Value = 134133.66
$ws.Range("I32").Value = 203924.31
$ws.Range("J32").Value = 17815.867
$ws.Range("K32").Value = 203924.31
$ws.Range("L32").Value = 17815.867
$ws.Range("M32").Value = -203637.31
$ws.Range("N32").Value = -18389.867
$ws.Range("H45").Value = 1938.6
$ws.Range("I45").Value = 1897.6666
$ws.Range("K45").Value = 1897.6666
$ws.Range("M45").Value = -1520.6666
$ws.Range("H61").Value = 1392259.6
$ws.Range("I61").Value = 3418.8
$ws.Range("J61").Value = 5885568
$ws.Range("K61").Value = 3418.8
$ws.Range("L61").Value = 5885568
$ws.Range("M61").Value = -3206.8
$ws.Range("N61").Value = -5885992
$ws.Range("H64").Value = 50000
$ws.Range("J64").Value = 50000
$ws.Range("L64").Value = 50000
$ws.Range("N64").Value = -50496
$ws.Range("H67").Value = 50000
$ws.Range("J67").Value = 50000
$ws.Range("L67").Value = 50000
$ws.Range("N67").Value = -51716
$ws.Range("H74").Value = 1924404.5
$ws.Range("I74").Value = 2529888.2
$ws.Range("J74").Value = 21456.143
$ws.Range("K74").Value = 2529888.2
$ws.Range("L74").Value = 21456.143
$ws.Range("M74").Value = -2529014.2
$ws.Range("N74").Value = -23204.143
$ws.Range("H77").Value = 1924404.5
$ws.Range("I77").Value = 2529888.2
$ws.Range("J77").Value = 21456.143
$ws.Range("K77").Value = 12649441
$ws.Range("L77").Value = 107280.715
$ws.Range("M77").Value = -12645073
$ws.Range("N77").Value = -116016.715
$ws.Range("H88").Value = 2237.92
$ws.Range("I88").Value = 1634.5
$ws.Range("J88").Value = 2640.2
$ws.Range("K88").Value = 1634.5
$ws.Range("L88").Value = 2640.2
$ws.Range("M88").Value = -1228.5
$ws.Range("N88").Value = -3452.2
$ws.Range("H91").Value = 2237.92
$ws.Range("I91").Value = 1634.5
$ws.Range("J91").Value = 2640.2
$ws.Range("K91").Value = 1634.5
$ws.Range("L91").Value = 2640.2
$ws.Range("M91").Value = -230.5
$ws.Range("N91").Value = -5448.2
$ws.Range("H97").Value = 25642046
$ws.Range("I97").Value = 794.5625
$ws.Range("J97").Value = 142859180
$ws.Range("K97").Value = 794.5625
$ws.Range("L97").Value = 142859180
$ws.Range("M97").Value = -298.5625
$ws.Range("N97").Value = -142860172
$ws.Range("H102").Value = 2057.8
$ws.Range("I102").Value = 2197.4614
$ws.Range("K102").Value = 2197.4614
$ws.Range("M102").Value = -575.4614000000001
$ws.Range("H103").Value = 20000
$ws.Range("J103").Value = 20000
$ws.Range("L103").Value = 20000
$ws.Range("N103").Value = -22344
$ws.Range("H110").Value = 1118.8667
$ws.Range("I110").Value = 839.75
$ws.Range("J110").Value = 1437.8572
$ws.Range("K110").Value = 839.75
$ws.Range("L110").Value = 1437.8572
$ws.Range("M110").Value = 1205.25
$ws.Range("N110").Value = -5527.8572
$ws.Range("H112").Value = 21997.4
$ws.Range("J112").Value = 21997.4
$ws.Range("L112").Value = 21997.4
$ws.Range("N112").Value = -24951.4
$ws.Range("H116").Value = 1146.6666
$ws.Range("I116").Value = 970
$ws.Range("J116").Value = 1500
$ws.Range("K116").Value = 970
$ws.Range("L116").Value = 1500
$ws.Range("M116").Value = 1324
$ws.Range("N116").Value = -6088
$ws.Range("H132").Value = 837161.8
$ws.Range("I132").Value = 1045743.44
$ws.Range("J132").Value = 2835.5
$ws.Range("K132").Value = 3137230.32
$ws.Range("L132").Value = 8506.5
$ws.Range("M132").Value = -3134700.32
$ws.Range("N132").Value = -13566.5
$ws.Range("H135").Value = 50429
$ws.Range("J135").Value = 50429
$ws.Range("L135").Value = 50429
$ws.Range("N135").Value = -60569
$ws.Range("H136").Value = 1392259.6
$ws.Range("I136").Value = 3418.8
$ws.Range("J136").Value = 5885568
$ws.Range("K136").Value = 10256.4
$ws.Range("L136").Value = 17656704
$ws.Range("M136").Value = -7706.400000000001
$ws.Range("N136").Value = -17661804

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1146.6666
$ws.Range("I3").Value = 970
$ws.Range("J3").Value = 1500
$ws.Range("K3").Value = 970
$ws.Range("L3").Value = 1500
$ws.Range("M3").Value = -856
$ws.Range("N3").Value = -1728
$ws.Range("H62").Value = 44173
$ws.Range("J62").Value = 44173
$ws.Range("L62").Value = 44173
$ws.Range("N62").Value = -45545
$ws.Range("H64").Value = 8757.666999999999
$ws.Range("J64").Value = 2509.5
$ws.Range("L64").Value = 2509.5
$ws.Range("N64").Value = -2959.5
$ws.Range("H65").Value = 44173
$ws.Range("J65").Value = 44173
$ws.Range("L65").Value = 132519
$ws.Range("N65").Value = -139383
$ws.Range("H67").Value = 8757.666999999999
$ws.Range("J67").Value = 2509.5
$ws.Range("L67").Value = 2509.5
$ws.Range("N67").Value = -4069.5
$ws.Range("H74").Value = 98500
$ws.Range("J74").Value = 98500
$ws.Range("L74").Value = 98500
$ws.Range("N74").Value = -100372
$ws.Range("H77").Value = 98500
$ws.Range("J77").Value = 98500
$ws.Range("L77").Value = 295500
$ws.Range("N77").Value = -304860
$ws.Range("H94").Value = 1659.4828
$ws.Range("I94").Value = 1670.6957
$ws.Range("J94").Value = 1616.5
$ws.Range("K94").Value = 1670.6957
$ws.Range("L94").Value = 1616.5
$ws.Range("M94").Value = -1219.6957
$ws.Range("N94").Value = -2518.5
$ws.Range("H99").Value = 6596.75
$ws.Range("I99").Value = 9202.691999999999
$ws.Range("J99").Value = 1757.1428
$ws.Range("K99").Value = 9202.691999999999
$ws.Range("L99").Value = 1757.1428
$ws.Range("M99").Value = -7704.691999999999
$ws.Range("N99").Value = -4753.1428
$ws.Range("H105").Value = 6488.5713
$ws.Range("I105").Value = 9750
$ws.Range("K105").Value = 9750
$ws.Range("M105").Value = -8003
$ws.Range("H124").Value = 40000
$ws.Range("J124").Value = 40000
$ws.Range("L124").Value = 40000
$ws.Range("N124").Value = -49820
$ws.Range("H134").Value = 4909783.5
$ws.Range("I134").Value = 5172.033
$ws.Range("J134").Value = 41694370
$ws.Range("K134").Value = 15516.099
$ws.Range("L134").Value = 125083110
$ws.Range("M134").Value = -12981.099
$ws.Range("N134").Value = -125088180

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H10").Value = 1176.5
$ws.Range("J10").Value = 5060
$ws.Range("L10").Value = 5060
$ws.Range("N10").Value = -5338
$ws.Range("H16").Value = 114028.89
$ws.Range("I16").Value = 3182
$ws.Range("J16").Value = 252587.5
$ws.Range("K16").Value = 3182
$ws.Range("L16").Value = 252587.5
$ws.Range("M16").Value = -2895
$ws.Range("N16").Value = -253161.5
$ws.Range("H22").Value = 1165.7
$ws.Range("I22").Value = 256.45456
$ws.Range("J22").Value = 2277
$ws.Range("K22").Value = 256.45456
$ws.Range("L22").Value = 2277
$ws.Range("M22").Value = 93.54543999999999
$ws.Range("N22").Value = -2977
$ws.Range("H28").Value = 52816.5
$ws.Range("J28").Value = 52816.5
$ws.Range("L28").Value = 52816.5
$ws.Range("N28").Value = -53306.5
$ws.Range("H31").Value = 1794812.8
$ws.Range("I31").Value = 1986542.6
$ws.Range("K31").Value = 1986542.6
$ws.Range("M31").Value = -1986247.6
$ws.Range("H34").Value = 1794812.8
$ws.Range("I34").Value = 1986542.6
$ws.Range("K34").Value = 1986542.6
$ws.Range("M34").Value = -1986340.6
$ws.Range("H43").Value = 16532.334
$ws.Range("I43").Value = 14597
$ws.Range("J43").Value = 17500
$ws.Range("K43").Value = 14597
$ws.Range("L43").Value = 17500
$ws.Range("M43").Value = -14413
$ws.Range("N43").Value = -17868
$ws.Range("H52").Value = 70999.5
$ws.Range("J52").Value = 99000
$ws.Range("L52").Value = 99000
$ws.Range("N52").Value = -99588
$ws.Range("H58").Value = 4913822
$ws.Range("I58").Value = 7098.4
$ws.Range("J58").Value = 6958290.5
$ws.Range("K58").Value = 7098.4
$ws.Range("L58").Value = 6958290.5
$ws.Range("M58").Value = -6895.4
$ws.Range("N58").Value = -6958696.5
$ws.Range("H59").Value = 49999
$ws.Range("I59").Value = 0
$ws.Range("K59").Value = 0
$ws.Range("M59").ClearContents()
$ws.Range("H62").Value = 3360.2727
$ws.Range("J62").Value = 3447.5
$ws.Range("L62").Value = 3447.5
$ws.Range("N62").Value = -4695.5
$ws.Range("H65").Value = 3360.2727
$ws.Range("J65").Value = 3447.5
$ws.Range("L65").Value = 17237.5
$ws.Range("N65").Value = -23477.5
$ws.Range("H87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("L87").ClearContents()
$ws.Range("N87").ClearContents()
$ws.Range("H90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("L90").ClearContents()
$ws.Range("N90").ClearContents()
$ws.Range("H93").Value = 11053.625
$ws.Range("I93").Value = 11053.625
$ws.Range("K93").Value = 11053.625
$ws.Range("M93").Value = -9181.625
$ws.Range("H95").Value = 32669.2
$ws.Range("J95").Value = 32669.2
$ws.Range("L95").Value = 32669.2
$ws.Range("N95").Value = -38161.2
$ws.Range("H99").Value = 5843037
$ws.Range("I99").Value = 25987.334
$ws.Range("J99").Value = 12387218
$ws.Range("K99").Value = 25987.334
$ws.Range("L99").Value = 12387218
$ws.Range("M99").Value = -24489.334
$ws.Range("N99").Value = -12390214
$ws.Range("H101").Value = 16532.334
$ws.Range("I101").Value = 14597
$ws.Range("J101").Value = 17500
$ws.Range("K101").Value = 14597
$ws.Range("L101").Value = 17500
$ws.Range("M101").Value = -11352
$ws.Range("N101").Value = -23990
$ws.Range("H105").Value = 7750.647
$ws.Range("I105").Value = 10697.728
$ws.Range("J105").Value = 2347.6667
$ws.Range("K105").Value = 10697.728
$ws.Range("L105").Value = 2347.6667
$ws.Range("M105").Value = -8950.727999999999
$ws.Range("N105").Value = -5841.6667
$ws.Range("H113").Value = 114028.89
$ws.Range("I113").Value = 3182
$ws.Range("J113").Value = 252587.5
$ws.Range("K113").Value = 3182
$ws.Range("L113").Value = 252587.5
$ws.Range("M113").Value = -1012
$ws.Range("N113").Value = -256927.5
$ws.Range("H126").Value = 5843037
$ws.Range("I126").Value = 25987.334
$ws.Range("J126").Value = 12387218
$ws.Range("K126").Value = 77962.00199999999
$ws.Range("L126").Value = 37161654
$ws.Range("M126").Value = -75492.00199999999
$ws.Range("N126").Value = -37166594
$ws.Range("H132").Value = 1266.6666
$ws.Range("I132").Value = 1266.6666
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 3799.9998
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -1269.9998
$ws.Range("N132").ClearContents()
$ws.Range("H134").Value = 1556.1305
$ws.Range("I134").Value = 1466.1904
$ws.Range("K134").Value = 4398.5712
$ws.Range("M134").Value = -1863.5712
$ws.Range("H136").Value = 4913822
$ws.Range("I136").Value = 7098.4
$ws.Range("J136").Value = 6958290.5
$ws.Range("K136").Value = 21295.2
$ws.Range("L136").Value = 20874871.5
$ws.Range("M136").Value = -18745.2
$ws.Range("N136").Value = -20879971.5

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 306.25
$ws.Range("I2").Value = 257.5
$ws.Range("K2").Value = 1545
$ws.Range("M2").Value = -1432
$ws.Range("H5").Value = 2185680.2
$ws.Range("I5").Value = 2551563.2
$ws.Range("J5").Value = 1791652.8
$ws.Range("K5").Value = 7654689.600000001
$ws.Range("L5").Value = 5374958.4
$ws.Range("M5").Value = -7654577.600000001
$ws.Range("N5").Value = -5375182.4
$ws.Range("H7").Value = 67.333336
$ws.Range("I7").Value = 67.333336
$ws.Range("K7").Value = 202.000008
$ws.Range("M7").Value = -90.00000800000001
$ws.Range("H10").Value = 635.7
$ws.Range("I10").Value = 150.25
$ws.Range("J10").Value = 959.3333
$ws.Range("K10").Value = 450.75
$ws.Range("L10").Value = 2877.9999
$ws.Range("M10").Value = -311.75
$ws.Range("N10").Value = -3155.9999
$ws.Range("H11").Value = 142857180
$ws.Range("J11").Value = 1000000000
$ws.Range("L11").Value = 3000000000
$ws.Range("N11").Value = -3000000280
$ws.Range("H38").Value = 56.285713
$ws.Range("I38").Value = 43.70968
$ws.Range("K38").Value = 131.12904
$ws.Range("M38").Value = 215.87096
$ws.Range("H47").Value = 204
$ws.Range("I47").Value = 204
$ws.Range("K47").Value = 612
$ws.Range("M47").Value = -181
$ws.Range("H105").Value = 18500
$ws.Range("J105").Value = 25000
$ws.Range("L105").Value = 75000
$ws.Range("N105").Value = -80242
$ws.Range("H107").Value = 4410.25
$ws.Range("I107").Value = 541.6667
$ws.Range("J107").Value = 5699.778
$ws.Range("K107").Value = 1625.0001
$ws.Range("L107").Value = 17099.334
$ws.Range("M107").Value = 294.9999
$ws.Range("N107").Value = -20939.334
$ws.Range("H111").Value = 7510.3335
$ws.Range("I111").Value = 854
$ws.Range("K111").Value = 2562
$ws.Range("M111").Value = 505
$ws.Range("H119").Value = 9479.637000000001
$ws.Range("I119").Value = 2379.3333
$ws.Range("K119").Value = 7137.999899999999
$ws.Range("M119").Value = -2299.999899999999
$ws.Range("H120").Value = 25362.9
$ws.Range("I120").Value = 23325.8
$ws.Range("J120").Value = 27400
$ws.Range("K120").Value = 69977.39999999999
$ws.Range("L120").Value = 82200
$ws.Range("M120").Value = -65139.39999999999
$ws.Range("N120").Value = -91876
$ws.Range("H121").Value = 55004548
$ws.Range("J121").Value = 55004548
$ws.Range("L121").Value = 165013644
$ws.Range("N121").Value = -165016264
$ws.Range("H124").Value = 11437.5
$ws.Range("I124").Value = 2000
$ws.Range("K124").Value = 6000
$ws.Range("M124").Value = -1090
$ws.Range("H129").Value = 6735607.5
$ws.Range("J129").Value = 9268698
$ws.Range("L129").Value = 27806094
$ws.Range("N129").Value = -27816094
$ws.Range("H131").Value = 5852337
$ws.Range("I131").Value = 1493.5555
$ws.Range("J131").Value = 11118096
$ws.Range("K131").Value = 4480.666499999999
$ws.Range("L131").Value = 33354288
$ws.Range("M131").Value = 559.3335000000006
$ws.Range("N131").Value = -33364368
$ws.Range("H134").Value = 5013
$ws.Range("I134").Value = 3418.4
$ws.Range("J134").Value = 8999.5
$ws.Range("K134").Value = 10255.2
$ws.Range("L134").Value = 26998.5
$ws.Range("M134").Value = -5185.200000000001
$ws.Range("N134").Value = -37138.5
$ws.Range("H135").Value = 2185680.2
$ws.Range("I135").Value = 2551563.2
$ws.Range("J135").Value = 1791652.8
$ws.Range("K135").Value = 22964068.8
$ws.Range("L135").Value = 16124875.2
$ws.Range("M135").Value = -22961533.8
$ws.Range("N135").Value = -16129945.2
$ws.Range("H136").Value = 7445.722
$ws.Range("I136").Value = 3668.75
$ws.Range("K136").Value = 11006.25
$ws.Range("M136").Value = -5906.25
$ws.Range("H137").Value = 6086.4517
$ws.Range("I137").Value = 2263.2
$ws.Range("K137").Value = 6789.599999999999
$ws.Range("M137").Value = -1689.599999999999
$ws.Range("H140").Value = 1499.0416
$ws.Range("I140").Value = 1117.9048
$ws.Range("K140").Value = 3353.7144
$ws.Range("M140").Value = 1826.2856

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H18").Value = 13898.5
$ws.Range("I18").Value = 13898.5
$ws.Range("K18").Value = 13898.5
$ws.Range("M18").Value = -13605.5
$ws.Range("H23").Value = 11814.667
$ws.Range("I23").Value = 10000
$ws.Range("J23").Value = 12333.143
$ws.Range("K23").Value = 10000
$ws.Range("L23").Value = 12333.143
$ws.Range("M23").Value = -9777
$ws.Range("N23").Value = -12779.143
$ws.Range("H70").Value = 67830.78
$ws.Range("I70").Value = 16184.8
$ws.Range("K70").Value = 16184.8
$ws.Range("M70").Value = -15914.8
$ws.Range("H73").Value = 67830.78
$ws.Range("I73").Value = 16184.8
$ws.Range("K73").Value = 16184.8
$ws.Range("M73").Value = -15248.8
$ws.Range("H80").Value = 2369.3333
$ws.Range("I80").Value = 2469.6667
$ws.Range("J80").Value = 2168.6667
$ws.Range("K80").Value = 2469.6667
$ws.Range("L80").Value = 2168.6667
$ws.Range("M80").Value = -1471.6667
$ws.Range("N80").Value = -4164.6667
$ws.Range("H83").Value = 2369.3333
$ws.Range("I83").Value = 2469.6667
$ws.Range("J83").Value = 2168.6667
$ws.Range("K83").Value = 12348.3335
$ws.Range("L83").Value = 10843.3335
$ws.Range("M83").Value = -7356.333500000001
$ws.Range("N83").Value = -20827.3335
$ws.Range("H97").Value = 1468.2433
$ws.Range("I97").Value = 1368.3448
$ws.Range("K97").Value = 1368.3448
$ws.Range("M97").Value = -872.3448000000001
$ws.Range("H113").Value = 3480.2307
$ws.Range("I113").Value = 3029.5293
$ws.Range("K113").Value = 3029.5293
$ws.Range("M113").Value = -859.5293000000001
$ws.Range("H122").Value = 6122.8076
$ws.Range("I122").Value = 6936.136
$ws.Range("J122").Value = 1649.5
$ws.Range("K122").Value = 20808.408
$ws.Range("L122").Value = 4948.5
$ws.Range("M122").Value = -18358.408
$ws.Range("N122").Value = -9848.5
$ws.Range("H126").Value = 6055.3335
$ws.Range("I126").Value = 6412.25
$ws.Range("J126").Value = 3200
$ws.Range("K126").Value = 19236.75
$ws.Range("L126").Value = 9600
$ws.Range("M126").Value = -16766.75
$ws.Range("N126").Value = -14540
$ws.Range("H132").Value = 16174.526
$ws.Range("I132").Value = 11517.556
$ws.Range("J132").Value = 100000
$ws.Range("K132").Value = 34552.66800000001
$ws.Range("L132").Value = 300000
$ws.Range("M132").Value = -32022.66800000001
$ws.Range("N132").Value = -305060
$ws.Range("H138").Value = 79997.5
$ws.Range("J138").Value = 79997.5
$ws.Range("L138").Value = 79997.5
$ws.Range("N138").Value = -90277.5
$ws.Range("H141").Value = 141332.67
$ws.Range("J141").Value = 141332.67
$ws.Range("L141").Value = 141332.67
$ws.Range("N141").Value = -151692.67

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 15138
$ws.Range("I7").Value = 15138
$ws.Range("K7").Value = 15138
$ws.Range("M7").Value = -15026
$ws.Range("H16").Value = 1535.5172
$ws.Range("I16").Value = 516.53845
$ws.Range("J16").Value = 10366.667
$ws.Range("K16").Value = 516.53845
$ws.Range("L16").Value = 10366.667
$ws.Range("M16").Value = -346.53845
$ws.Range("N16").Value = -10706.667
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 0
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 0
$ws.Range("M17").ClearContents()
$ws.Range("N17").ClearContents()
$ws.Range("H38").Value = 30000
$ws.Range("J38").Value = 30000
$ws.Range("L38").Value = 30000
$ws.Range("N38").Value = -30820
$ws.Range("H40").Value = 6414.143
$ws.Range("I40").Value = 6180.8
$ws.Range("K40").Value = 6180.8
$ws.Range("M40").Value = -6044.8
$ws.Range("H46").Value = 2410.5356
$ws.Range("I46").Value = 350.93332
$ws.Range("J46").Value = 4787
$ws.Range("K46").Value = 350.93332
$ws.Range("L46").Value = 4787
$ws.Range("M46").Value = -162.93332
$ws.Range("N46").Value = -5163
$ws.Range("H55").Value = 1673.1904
$ws.Range("I55").Value = 1911.2222
$ws.Range("J55").Value = 1494.6666
$ws.Range("K55").Value = 1911.2222
$ws.Range("L55").Value = 1494.6666
$ws.Range("M55").Value = -1738.2222
$ws.Range("N55").Value = -1840.6666
$ws.Range("H61").Value = 12217.588
$ws.Range("I61").Value = 11817.929
$ws.Range("J61").Value = 14082.667
$ws.Range("K61").Value = 11817.929
$ws.Range("L61").Value = 14082.667
$ws.Range("M61").Value = -11615.929
$ws.Range("N61").Value = -14486.667
$ws.Range("H68").Value = 2899.8
$ws.Range("I68").Value = 1999.75
$ws.Range("K68").Value = 1999.75
$ws.Range("M68").Value = -1250.75
$ws.Range("H71").Value = 2899.8
$ws.Range("I71").Value = 1999.75
$ws.Range("K71").Value = 9998.75
$ws.Range("M71").Value = -6254.75
$ws.Range("H82").Value = 2709.75
$ws.Range("I82").Value = 3090.5715
$ws.Range("J82").Value = 1821.1666
$ws.Range("K82").Value = 3090.5715
$ws.Range("L82").Value = 1821.1666
$ws.Range("M82").Value = -2729.5715
$ws.Range("N82").Value = -2543.1666
$ws.Range("H85").Value = 2709.75
$ws.Range("I85").Value = 3090.5715
$ws.Range("J85").Value = 1821.1666
$ws.Range("K85").Value = 3090.5715
$ws.Range("L85").Value = 1821.1666
$ws.Range("M85").Value = -1842.5715
$ws.Range("N85").Value = -4317.1666
$ws.Range("H93").Value = 1592.1305
$ws.Range("I93").Value = 1273.9412
$ws.Range("J93").Value = 2493.6667
$ws.Range("K93").Value = 1273.9412
$ws.Range("L93").Value = 2493.6667
$ws.Range("M93").Value = -25.94119999999998
$ws.Range("N93").Value = -4989.6667
$ws.Range("H100").Value = 1476.091
$ws.Range("I100").Value = 1407
$ws.Range("K100").Value = 1407
$ws.Range("M100").Value = -866
$ws.Range("H113").Value = 12217.588
$ws.Range("I113").Value = 11817.929
$ws.Range("J113").Value = 14082.667
$ws.Range("K113").Value = 11817.929
$ws.Range("L113").Value = 14082.667
$ws.Range("M113").Value = -9647.929
$ws.Range("N113").Value = -18422.667
$ws.Range("H122").Value = 3159.3
$ws.Range("I122").Value = 2451
$ws.Range("K122").Value = 7353
$ws.Range("M122").Value = -4903
$ws.Range("H126").Value = 15138
$ws.Range("I126").Value = 15138
$ws.Range("K126").Value = 45414
$ws.Range("M126").Value = -42944
$ws.Range("H132").Value = 6498477.5
$ws.Range("I132").Value = 11687083
$ws.Range("J132").Value = 12720.875
$ws.Range("K132").Value = 35061249
$ws.Range("L132").Value = 38162.625
$ws.Range("M132").Value = -35058719
$ws.Range("N132").Value = -43222.625
$ws.Range("H133").Value = 78826.375
$ws.Range("J133").Value = 71994.75
$ws.Range("L133").Value = 71994.75
$ws.Range("N133").Value = -77054.75
$ws.Range("H136").Value = 5324719
$ws.Range("I136").Value = 5685385
$ws.Range("K136").Value = 17056155
$ws.Range("M136").Value = -17053605

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 363724.1
$ws.Range("I4").Value = 444518.66
$ws.Range("K4").Value = 444518.66
$ws.Range("M4").Value = -444405.66
$ws.Range("H22").Value = 0
$ws.Range("J22").Value = 0
$ws.Range("L22").ClearContents()
$ws.Range("N22").ClearContents()
$ws.Range("H53").Value = 50
$ws.Range("I53").Value = 50
$ws.Range("K53").Value = 50
$ws.Range("M53").Value = 557
$ws.Range("H62").Value = 6993.3125
$ws.Range("I62").Value = 5586.8
$ws.Range("J62").Value = 7632.636
$ws.Range("K62").Value = 5586.8
$ws.Range("L62").Value = 7632.636
$ws.Range("M62").Value = -4962.8
$ws.Range("N62").Value = -8880.636
$ws.Range("H65").Value = 6993.3125
$ws.Range("I65").Value = 5586.8
$ws.Range("J65").Value = 7632.636
$ws.Range("K65").Value = 27934
$ws.Range("L65").Value = 38163.18
$ws.Range("M65").Value = -24814
$ws.Range("N65").Value = -44403.18
$ws.Range("H81").Value = 4721.4443
$ws.Range("I81").Value = 1480
$ws.Range("K81").Value = 2960
$ws.Range("M81").Value = -1899
$ws.Range("H84").Value = 4721.4443
$ws.Range("I84").Value = 1480
$ws.Range("K84").Value = 14800
$ws.Range("M84").Value = -9496
$ws.Range("H100").Value = 1705.5
$ws.Range("I100").Value = 1705.5
$ws.Range("K100").Value = 3411
$ws.Range("M100").Value = -2870
$ws.Range("H107").Value = 2040
$ws.Range("I107").Value = 1184.762
$ws.Range("K107").Value = 3554.286
$ws.Range("M107").Value = -1634.286
$ws.Range("H113").Value = 785.97675
$ws.Range("I113").Value = 786.3103599999999
$ws.Range("J113").Value = 785.2857
$ws.Range("K113").Value = 2358.93108
$ws.Range("L113").Value = 2355.8571
$ws.Range("M113").Value = -188.9310799999998
$ws.Range("N113").Value = -6695.8571
$ws.Range("H122").Value = 57764.15
$ws.Range("I122").Value = 1364.6
$ws.Range("K122").Value = 4093.8
$ws.Range("M122").Value = -1643.8
$ws.Range("H123").Value = 51449.5
$ws.Range("J123").Value = 51449.5
$ws.Range("L123").Value = 51449.5
$ws.Range("N123").Value = -61249.5
$ws.Range("H132").Value = 6668753.5
$ws.Range("I132").Value = 7938301.5
$ws.Range("J132").Value = 3626
$ws.Range("K132").Value = 23814904.5
$ws.Range("L132").Value = 10878
$ws.Range("M132").Value = -23812374.5
$ws.Range("N132").Value = -15938
$ws.Range("H136").Value = 5293912
$ws.Range("I136").Value = 1175852.1
$ws.Range("J136").Value = 22223714
$ws.Range("K136").Value = 3527556.3
$ws.Range("L136").Value = 66671142
$ws.Range("M136").Value = -3525006.3
$ws.Range("N136").Value = -66676242
$ws.Range("H140").Value = 88586
$ws.Range("J140").Value = 88586
$ws.Range("L140").Value = 88586
$ws.Range("N140").Value = -98946
$ws.Range("H141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("L141").ClearContents()
$ws.Range("N141").ClearContents()

Write-Output "Applied changes: sets and clears complete"